$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of order data to append (rows 18-21).
# All values in this sheet are stored as text (inline strings), including
# numeric-looking quantities/prices, so force text number format before
# writing the values to avoid Excel auto-converting them to numbers.
$newRows = @(
    @("TN374", "Natalie's - Lemonade", "1", "9.30", "9.30"),
    @("AH252", "Natalie's - Orange Juice", "2", "24.50", "49.00"),
    @("TN454", "Natalie's - Orange Mango", "1", "13.38", "13.38"),
    @("TN380", "Natalie's - Strawberry Lemonade", "1", "10.15", "10.15")
)

$startRow = 18
$endRow = $startRow + $newRows.Length - 1

# Force the target range to text format so numeric-looking strings are
# preserved as text rather than being coerced into numbers.
$ws.Range("A$startRow`:E$endRow").NumberFormat = "@"

$r = $startRow
foreach ($rowData in $newRows) {
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $r = $r + 1
}
